# Update "想去人数" (want-to-go count, column F) values across all four
# worksheets to reflect the freshly scraped numbers (gh-pages rebuild at
# commit 456a3b4). Only column F numeric values change; everything else
# (labels, G/price column, styles, etc.) stays untouched.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$Updates = @{
    4  = 1295
    6  = 355
    7  = 1163
    8  = 442
    9  = 7092
    12 = 2041
    13 = 7963
    16 = 5503
    17 = 49
    18 = 2398
    19 = 1021
    20 = 4561
    25 = 367
    26 = 254
    28 = 2318
    29 = 24
    30 = 262
    31 = 75
    32 = 136
    33 = 575
    34 = 4
    36 = 1484
    37 = 31
    38 = 4
    39 = 2301
    40 = 2210
    42 = 5
}
foreach ($row in $Updates.Keys) {
    $ws.Range("F$row").Value = $Updates[$row]
}

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$Updates = @{
    3 = 74
    4 = 59
    8 = 96
}
foreach ($row in $Updates.Keys) {
    $ws.Range("F$row").Value = $Updates[$row]
}

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$Updates = @{
    3 = 1277
}
foreach ($row in $Updates.Keys) {
    $ws.Range("F$row").Value = $Updates[$row]
}

# Sheet "全部类型" (All types combined)
$ws = $wb.Worksheets.Item("全部类型")
$Updates = @{
    4  = 1277
    5  = 1295
    7  = 355
    8  = 1163
    9  = 442
    10 = 7092
    13 = 2041
    14 = 7963
    17 = 5503
    18 = 49
    19 = 2398
    20 = 1021
    21 = 4561
    25 = 74
    27 = 59
    28 = 367
    29 = 254
    30 = 2318
    31 = 24
    32 = 262
    33 = 75
    34 = 136
    36 = 575
    37 = 4
    40 = 1484
    41 = 31
    42 = 4
    43 = 2301
    45 = 2210
    47 = 5
    49 = 96
}
foreach ($row in $Updates.Keys) {
    $ws.Range("F$row").Value = $Updates[$row]
}

$wb.Save()
